$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure text-like numeric strings (prices) are written verbatim, not coerced to floats.
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '29.728.69'
$ws.Range('E2').Value = '  -1.12%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.103.45'
$ws.Range('E3').Value = '  -0.11%  '
$ws.Range('E4').Value = '  +0.53%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '346.75'
$ws.Range('E5').Value = '  -0.35%  '
$ws.Range('E6').Value = '  +0.52%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5192'
$ws.Range('E7').Value = '  +0.43%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.4449'
$ws.Range('E8').Value = '  +0.10%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '54.32'
$ws.Range('E9').Value = '  +3.96%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.09335'
$ws.Range('E10').Value = '  +3.99%  '
$ws.Range('E11').Value = '  +0.34%  '
$ws.Range('E12').Value = '  -2.33%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '2.156.37'
$ws.Range('E13').Value = '  +2.33%  '
$ws.Range('E14').Value = '  +1.29%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '8.296'
$ws.Range('E15').Value = '  +0.60%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '102.49'
$ws.Range('E16').Value = '  +3.25%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.00001158'
$ws.Range('E17').Value = '  +0.64%  '
$ws.Range('E18').Value = '  +0.46%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '21.31'
$ws.Range('E19').Value = '  +1.90%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.06677'
$ws.Range('E20').Value = '  +0.00%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.298'
$ws.Range('E21').Value = '  +0.88%  '
$ws.Range('E22').Value = '  +0.42%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '29.789.57'
$ws.Range('E24').Value = '  -0.59%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.327'
$ws.Range('E25').Value = '  -0.92%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.395.34'
$ws.Range('E26').Value = '  +1.66%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '22.04'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.541'
$ws.Range('E28').Value = '  -0.25%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '162.25'
$ws.Range('E29').Value = '  -0.14%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '133.77'
$ws.Range('E30').Value = '  +0.12%  '
$ws.Range('B31').Value = 'ARBITRUM'
$ws.Range('C31').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.796'
$ws.Range('E31').Value = '  +9.40%  '
$ws.Range('B32').Value = 'ImmutableX'
$ws.Range('C32').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.145'
$ws.Range('E32').Value = '  -2.68%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.1055'
$ws.Range('E33').Value = '  -1.07%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '6.225'
$ws.Range('E34').Value = '  -0.20%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.947'
$ws.Range('E35').Value = '  -0.30%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '6.334'
$ws.Range('E36').Value = '  +7.06%  '
$ws.Range('E37').Value = '  +5.33%  '
$ws.Range('E38').Value = '  +0.59%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.06770'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.7021'
$ws.Range('E40').Value = '  +2.93%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '12.58'
$ws.Range('E41').Value = '  +0.18%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.331'
$ws.Range('E42').Value = '  +2.94%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.2230'
$ws.Range('E43').Value = '  -3.11%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.6824'
$ws.Range('E44').Value = '  +6.86%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '14.53'
$ws.Range('E45').Value = '  +1.70%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.357'
$ws.Range('E46').Value = '  +2.67%  '
$ws.Range('B47').Value = 'Frax'
$ws.Range('C47').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.006'
$ws.Range('E47').Value = '  +0.51%  '
$ws.Range('B48').Value = 'PancakeSwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '3.635'
$ws.Range('E48').Value = '  -0.12%  '
$ws.Range('B49').Value = 'BabyDogeCoin'
$ws.Range('C49').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.00000000356'
$ws.Range('E49').Value = '  -2.07%  '
$ws.Range('B50').Value = 'WEMIXTOKEN'
$ws.Range('C50').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('E50').Value = '  +4.93%  '
$ws.Range('B51').Value = 'EOS'
$ws.Range('C51').Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.221'
$ws.Range('E51').Value = '  -0.01%  '
